$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..10) {
    $ws.Range("B${r}:D${r}").Value = "N/A"
    $ws.Range("B${r}:D${r}").NumberFormat = "0.00"

    $ws.Range("E${r}").Value = "N/A"
    $ws.Range("E${r}").NumberFormat = "0%"

    $ws.Range("F${r}:G${r}").Value = "N/A"
    $ws.Range("F${r}:G${r}").NumberFormat = "@"
}

[void]$ws.Range("H15").Select()
